$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.412.52"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.045.63"
$ws.Range("E3").Value = "  +3.96%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.66"
$ws.Range("E5").Value = "  -1.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "625.00"
$ws.Range("E6").Value = "  +4.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("E9").Value = "  +4.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.039.96"
$ws.Range("E10").Value = "  +3.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.443"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("E13").Value = "  +5.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.604.76"
$ws.Range("E14").Value = "  +3.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.14"
$ws.Range("E15").Value = "  +3.80%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.343.35"
$ws.Range("E16").Value = "  +0.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000194"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.048.96"
$ws.Range("E18").Value = "  +4.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.61"
$ws.Range("E19").Value = "  +3.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.00"
$ws.Range("E20").Value = "  +1.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "375.55"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.32"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.38"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.26"
$ws.Range("E24").Value = "  +2.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.38"
$ws.Range("E27").Value = "  +2.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("E28").Value = "  +1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000110"
$ws.Range("E29").Value = "  +1.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.996"
$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +6.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  +1.32%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.95"
$ws.Range("E33").Value = "  +6.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "491.91"
$ws.Range("E34").Value = "  -2.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "20.69"
$ws.Range("E36").Value = "  +2.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.78"
$ws.Range("E37").Value = "  -0.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.04"
$ws.Range("E38").Value = "  +2.10%  "

$ws.Range("E39").Value = "  +3.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "190.99"
$ws.Range("E41").Value = "  +5.30%  "

$ws.Range("E42").Value = "  -6.57%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.800"
$ws.Range("E44").Value = "  +21.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.14"
$ws.Range("E45").Value = "  +3.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.27"
$ws.Range("E46").Value = "  +5.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "42.03"
$ws.Range("E47").Value = "  +4.65%  "

$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.47"
$ws.Range("E49").Value = "  +5.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.607"
$ws.Range("E50").Value = "  +4.49%  "

$ws.Range("E51").Value = "  +4.41%  "
